# Generate Report for Archive
#
# The localization job moved from "handoff" to "in translation", so the
# Status column on every sheet (the Overview roll-up's per-locale columns
# plus each locale sheet's own Status column) needs to reflect the new
# state, and the Status column should be re-sized to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview: column E = zh-cn status, column F = de-de status
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Locale sheets: column C = Status
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Resize the Status column(s) to fit the new (shorter) text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
